$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.195906432748538
$ws.Range("C2").Value = 0.5175438596491229
$ws.Range("J2").Value = 0.02631578947368421
$ws.Range("S2").Value = 0.1023391812865497
$ws.Range("C3").Value = 0.02197802197802198
$ws.Range("J3").Value = 0.03846153846153846
$ws.Range("P3").Value = 0.7307692307692307
$ws.Range("S3").Value = 0.2087912087912088
$ws.Range("B6").Value = 0.09051724137931035
$ws.Range("D6").Value = 0.01293103448275862
$ws.Range("F6").Value = 0.08620689655172414
$ws.Range("J6").Value = 0.1767241379310345
$ws.Range("O6").Value = 0.03448275862068965
$ws.Range("Q6").Value = 0.1293103448275862
$ws.Range("R6").Value = 0.04310344827586207
$ws.Range("S6").Value = 0.4267241379310345
$ws.Range("B7").Value = 0.107843137254902
$ws.Range("D7").Value = 0.01470588235294118
$ws.Range("F7").Value = 0.06372549019607843
$ws.Range("J7").Value = 0.107843137254902
$ws.Range("O7").Value = 0.004901960784313725
$ws.Range("Q7").Value = 0.09803921568627451
$ws.Range("S7").Value = 0.5196078431372549
$ws.Range("B8").Value = 0.134453781512605
$ws.Range("D8").Value = 0.008403361344537815
$ws.Range("E8").Value = 0.005602240896358543
$ws.Range("F8").Value = 0.06722689075630252
$ws.Range("J8").Value = 0.08403361344537816
$ws.Range("O8").Value = 0.02240896358543417
$ws.Range("Q8").Value = 0.2100840336134454
$ws.Range("R8").Value = 0.06722689075630252
$ws.Range("S8").Value = 0.4005602240896359
$ws.Range("B9").Value = 0.1216216216216216
$ws.Range("D9").Value = 0.02027027027027027
$ws.Range("F9").Value = 0.07432432432432433
$ws.Range("J9").Value = 0.1283783783783784
$ws.Range("O9").Value = 0.03378378378378379
$ws.Range("Q9").Value = 0.1959459459459459
$ws.Range("R9").Value = 0.07432432432432433
$ws.Range("S9").Value = 0.3513513513513514
$ws.Range("B10").Value = 0.1386735572782084
$ws.Range("D10").Value = 0.020671834625323
$ws.Range("E10").Value = 0.0008613264427217916
$ws.Range("F10").Value = 0.08871662360034453
$ws.Range("J10").Value = 0.1223083548664944
$ws.Range("O10").Value = 0.01636520241171404
$ws.Range("Q10").Value = 0.1739879414298019
$ws.Range("R10").Value = 0.07235142118863049
$ws.Range("S10").Value = 0.3660637381567614
$ws.Range("G11").Value = 0.1559633027522936
$ws.Range("J11").Value = 0.08868501529051988
$ws.Range("K11").Value = 0.1957186544342508
$ws.Range("L11").Value = 0.5412844036697247
$ws.Range("S11").Value = 0.01834862385321101
$ws.Range("G12").Value = 0.7307692307692307
$ws.Range("J12").Value = 0.1978021978021978
$ws.Range("K12").Value = 0.02197802197802198
$ws.Range("L12").Value = 0.01648351648351648
$ws.Range("S12").Value = 0.03296703296703297
$ws.Range("G13").Value = 0.6444444444444445
$ws.Range("J13").Value = 0.2888888888888889
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01777777777777778
$ws.Range("H15").Value = 0.1022222222222222
$ws.Range("I15").Value = 0.08
$ws.Range("J15").Value = 0.3688888888888889
$ws.Range("K15").Value = 0.07111111111111111
$ws.Range("M15").Value = 0.02222222222222222
$ws.Range("O15").Value = 0.08
$ws.Range("S15").Value = 0.2577777777777778
$ws.Range("F16").Value = 0.009615384615384616
$ws.Range("H16").Value = 0.1586538461538461
$ws.Range("I16").Value = 0.0576923076923077
$ws.Range("J16").Value = 0.4278846153846154
$ws.Range("K16").Value = 0.1346153846153846
$ws.Range("M16").Value = 0.004807692307692308
$ws.Range("N16").Value = 0.004807692307692308
$ws.Range("O16").Value = 0.0673076923076923
$ws.Range("S16").Value = 0.1346153846153846
$ws.Range("F17").Value = 0.0254957507082153
$ws.Range("H17").Value = 0.1218130311614731
$ws.Range("I17").Value = 0.08781869688385269
$ws.Range("J17").Value = 0.4192634560906516
$ws.Range("K17").Value = 0.1218130311614731
$ws.Range("M17").Value = 0.0169971671388102
$ws.Range("O17").Value = 0.05099150141643059
$ws.Range("S17").Value = 0.1558073654390935
$ws.Range("F18").Value = 0.006756756756756757
$ws.Range("H18").Value = 0.1216216216216216
$ws.Range("I18").Value = 0.0945945945945946
$ws.Range("J18").Value = 0.4459459459459459
$ws.Range("K18").Value = 0.1418918918918919
$ws.Range("M18").Value = 0.01351351351351351
$ws.Range("O18").Value = 0.06756756756756757
$ws.Range("S18").Value = 0.1081081081081081
$ws.Range("F19").Value = 0.01367658889782784
$ws.Range("H19").Value = 0.1971037811745776
$ws.Range("I19").Value = 0.0587288817377313
$ws.Range("J19").Value = 0.3547868061142397
$ws.Range("K19").Value = 0.1190667739340306
$ws.Range("M19").Value = 0.02574416733708769
$ws.Range("O19").Value = 0.08447304907481899
$ws.Range("S19").Value = 0.1464199517296862
